# Edit: remove the "Features" (Numbered Colored Cards / Action Cards) slide
# from the deck, and tweak the player-count text on the "Rules" slide from
# "2 players" to "2-10 players".

$p = $ppt.ActivePresentation

# --- 1. Delete the slide with id=280 ("Features") -------------------------
# It is the 4th slide in the deck (before any edits), right after "Rules".
$targetIndex = 0
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    if ($p.Slides.Item($i).SlideID -eq 280) {
        $targetIndex = $i
        break
    }
}
if ($targetIndex -gt 0) {
    $p.Slides.Item($targetIndex).Delete()
}

# --- 2. Update "2 players" -> "2-10 players" on the "Rules" slide ---------
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -like "*2 players*") {
                $paraCount = $tr.Paragraphs().Count
                for ($k = 1; $k -le $paraCount; $k++) {
                    $para = $tr.Paragraphs($k, 1)
                    if ($para.Text.TrimEnd() -eq "2 players") {
                        $lead = $para.Characters(1, 2)
                        $lead.Text = "2-10 "
                    }
                }
            }
        }
    }
}
